# Actualización automática 2025-11-12 11:30:07
# Incrementa la venta de PORCELANATO de noviembre para
# HIDALGO HIDALGO PEDRO GUSTAVO / MEGAMAFERS S.A. y propaga
# el cambio a los totales dependientes en las 3 hojas.

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" (sheet1): detalle por categoría ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M13").Value = 4857.52

# --- Hoja "VENTA MENSUAL" (sheet2): detalle mensual ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F13").Value = 4857.52
$wsMensual.Range("F23").Value = 7181.500000000001

# --- Hoja "CUMPLIMIENTO MENSUAL" (sheet3): resumen de cumplimiento ---
$wsCump = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCump.Range("D12").Value = 6818.45
$wsCump.Range("E12").Value = 37599.55
$wsCump.Range("F12").Value = 0.1535064613445

$wsCump.Range("D14").Value = 7181.5
$wsCump.Range("E14").Value = 48217.97101170095
$wsCump.Range("F14").Value = 0.12963120168572
